$d = $word.ActiveDocument

$replacements = @(
    @("600÷5=", "139÷6="),
    @("755÷8=", "538÷4="),
    @("692÷4=", "200÷3="),
    @("883÷8=", "117÷8="),
    @("957÷3=", "477÷3="),
    @("157÷7=", "443÷7="),
    @("829÷3=", "282÷9="),
    @("685÷9=", "958÷8="),
    @("601÷3=", "564÷6="),
    @("737÷9=", "860÷3="),
    @("940÷3=", "288÷6="),
    @("323÷2=", "193÷6="),
    @("905÷2=", "258÷7="),
    @("599÷6=", "814÷9="),
    @("580÷3=", "470÷3="),
    @("753÷5=", "388÷2="),
    @("474÷8=", "113÷4="),
    @("533÷6=", "397÷6="),
    @("324÷9=", "675÷4="),
    @("172÷2=", "905÷4="),
    @("776÷7=", "937÷6="),
    @("433÷6=", "216÷4="),
    @("606÷9=", "978÷5="),
    @("197÷9=", "805÷8="),
    @("453÷9=", "910÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
